$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: (empty) -> false (text "false", not boolean -- round-trip
# through a formula + paste-values so it lands as a text cell like the
# original data, instead of Excel's normal TRUE/FALSE auto-typing)
$cExperimental = $ws.Range("B7")
$cExperimental.Formula = "=""false"""
$cExperimental.Copy()
$cExperimental.PasteSpecial(-4163)

# Date: updated timestamp
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# Case Sensitive: (empty) -> true (text "true", see note above)
$cCaseSensitive = $ws.Range("B15")
$cCaseSensitive.Formula = "=""true"""
$cCaseSensitive.Copy()
$cCaseSensitive.PasteSpecial(-4163)
